# "criando classe para comandos sql"
#
# Duplicate the last "Classe: Database" slide (slide 13) to create a new
# "Classe: ComandoSQL" slide, landing it right before the final
# "Fluxograma da aplicação" slide (slide 14), which is exactly what
# PowerPoint's Slide.Duplicate() does: the copy is inserted immediately
# after its source, pushing everything after it one slot later.

$p = $ppt.ActivePresentation

$srcSlide = $p.Slides.Item(13)

$dupRange = $srcSlide.Duplicate()
$newSlide = $dupRange.Item(1)

# --- Shape 1: title bar ("Retângulo 3") ---------------------------------
# Was: "Classe" + ": Database"
# Now: "Classe" + ": " + "ComandoSQL"  (3 runs)
$titleShape = $newSlide.Shapes.Item(1)
$titleRange = $titleShape.TextFrame.TextRange
$fullLen = $titleRange.Length
# "Classe: " is 8 characters; replace everything after it.
$tail = $titleRange.Characters(9, $fullLen - 8)
$tail.Text = "ComandoSQL"

# --- Shape 4: big name inside the rounded rectangle ("CaixaDeTexto 4") --
$nameShape = $newSlide.Shapes.Item(4)
$nameShape.TextFrame.TextRange.Text = "ComandoSQL"

# --- Shape 5: properties bullet list ("CaixaDeTexto 5") -----------------
$propsShape = $newSlide.Shapes.Item(5)
$propsShape.TextFrame.TextRange.Text = "ComandoSQL"

# --- Shape 6: methods bullet list ("CaixaDeTexto 6") ---------------------
$methodsShape = $newSlide.Shapes.Item(6)
$methodsShape.TextFrame.TextRange.Text = ""
